$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 28 (2025Q2) with refreshed recurrence metrics
$ws.Range("C28").Value = 394
$ws.Range("D28").Value = 44
$ws.Range("E28").Value = 350
$ws.Range("F28").Value = 6.853582554517133
